$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" (copy H1's header formatting) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Body rows: I and J values for each data row (r => [I, J]) ---
$values = @{
    "2"  = @(1,5)
    "3"  = @(1,6)
    "4"  = @(1,5)
    "5"  = @(1,6)
    "6"  = @(1,8)
    "7"  = @(1,7)
    "8"  = @(1,7)
    "9"  = @(1,6)
    "10" = @(1,5)
    "11" = @(1,6)
    "12" = @(1,7)
    "13" = @(1,7)
    "14" = @(1,8)
    "15" = @(1,6)
    "16" = @(1,7)
    "17" = @(1,6)
    "18" = @(1,8)
    "19" = @(1,6)
    "20" = @(1,7)
    "21" = @(1,6)
    "22" = @(1,5)
    "23" = @(1,6)
    "24" = @(1,6)
    "25" = @(1,7)
    "26" = @(1,8)
    "27" = @(1,7)
    "28" = @(1,7)
    "29" = @(1,6)
    "30" = @(1,7)
    "31" = @(1,6)
    "32" = @(1,5)
    "33" = @(1,4)
    "34" = @(6,8)
    "35" = @(5,7)
    "36" = @(5,8)
    "37" = @(4,5)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item([int]$r, 9).Value = $pair[0]
    $ws.Cells.Item([int]$r, 10).Value = $pair[1]
}
